# #5: property boat&car done
# Adds a "capacity" column to the 汽車 (Car) sheet (new column C, pushing
# owner/register_date/register_reason/acquire_value one slot right) and
# fills in the extra identifying metadata columns (property_category,
# category, date, legislator_name, legislator_id, source_file, index)
# that the other sheets (e.g. 土地) already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# Insert a new column before column C; existing C..G (owner, register_date,
# register_reason, acquire_value) shift right to D..H, keeping their styles.
$ws.Columns.Item(3).Insert()

# --- Header row (row 1) ---
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "capacity"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "register_date"
$ws.Cells.Item(1, 6).Value = "register_reason"
$ws.Cells.Item(1, 7).Value = "acquire_value"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# --- Data row (row 2) ---
$ws.Cells.Item(2, 2).Value = "TOYOTA"
$ws.Cells.Item(2, 3).Value = 2362
$ws.Cells.Item(2, 4).Value = "陳歐珀"
$ws.Cells.Item(2, 5).Value = "101年02月"
$ws.Cells.Item(2, 6).Value = "買賣"
$ws.Cells.Item(2, 7).Value = 1500000
$ws.Cells.Item(2, 8).Value = "land"
$ws.Cells.Item(2, 9).Value = "normal"

# This column holds text that otherwise auto-parses as a date, so force
# text format before assigning.
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2012-05-01"

$ws.Cells.Item(2, 11).Value = "陳歐珀"
$ws.Cells.Item(2, 12).Value = 1753
$ws.Cells.Item(2, 13).Value = "tmpe4f31"
$ws.Cells.Item(2, 14).Value = 30

# The columns beyond the original A1:G2 range (H..N) don't inherit the
# sheet's header/data styling (bold+border vs. plain) automatically, so
# copy it over explicitly to match columns B..G.
for ($col = 8; $col -le 14; $col++) {
    $headerCell = $ws.Cells.Item(1, $col)
    $headerCell.Font.Bold = $true
    $headerCell.HorizontalAlignment = -4108  # xlCenter
    $headerCell.VerticalAlignment = -4160    # xlTop
    $headerCell.Borders.LineStyle = 1        # xlContinuous
}
